$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("144:147").Insert()

# Row 144
$ws.Range("A144").Value = 12
$ws.Range("B144").Value = 'Mapocho Venta Directa de Santiago'
$ws.Range("C144").Value = 'Metropolitana'
$ws.Range("D144").Value = 44637
$ws.Range("E144").Value = 13
$ws.Range("F144").Value = 100112043
$ws.Range("G144").Value = 'Pepino dulce'
$ws.Range("H144").Value = 'Cultivar IV Región'
$ws.Range("I144").Value = 'Especial'
$ws.Range("J144").Value = 250
$ws.Range("K144").Value = 15000
$ws.Range("L144").Value = 15000
$ws.Range("M144").Value = 15000
$ws.Range("N144").Value = '$/bandeja 18 kilos'
$ws.Range("O144").Value = 'Provincia de Limarí'
$ws.Range("P144").Value = 833
$ws.Range("Q144").Value = 18
$ws.Range("R144").Value = 'Hortaliza'

# Row 145
$ws.Range("A145").Value = 12
$ws.Range("B145").Value = 'Mapocho Venta Directa de Santiago'
$ws.Range("C145").Value = 'Metropolitana'
$ws.Range("D145").Value = 44637
$ws.Range("E145").Value = 13
$ws.Range("F145").Value = 100112043
$ws.Range("G145").Value = 'Pepino dulce'
$ws.Range("H145").Value = 'Cultivar IV Región'
$ws.Range("I145").Value = 'Primera'
$ws.Range("J145").Value = 280
$ws.Range("K145").Value = 13000
$ws.Range("L145").Value = 13000
$ws.Range("M145").Value = 13000
$ws.Range("N145").Value = '$/bandeja 18 kilos'
$ws.Range("O145").Value = 'Provincia de Limarí'
$ws.Range("P145").Value = 722
$ws.Range("Q145").Value = 18
$ws.Range("R145").Value = 'Hortaliza'

# Row 146
$ws.Range("A146").Value = 12
$ws.Range("B146").Value = 'Mapocho Venta Directa de Santiago'
$ws.Range("C146").Value = 'Metropolitana'
$ws.Range("D146").Value = 44637
$ws.Range("E146").Value = 13
$ws.Range("F146").Value = 100112043
$ws.Range("G146").Value = 'Pepino dulce'
$ws.Range("H146").Value = 'Cultivar IV Región'
$ws.Range("I146").Value = 'Segunda'
$ws.Range("J146").Value = 220
$ws.Range("K146").Value = 10000
$ws.Range("L146").Value = 10000
$ws.Range("M146").Value = 10000
$ws.Range("N146").Value = '$/bandeja 18 kilos'
$ws.Range("O146").Value = 'Provincia de Limarí'
$ws.Range("P146").Value = 556
$ws.Range("Q146").Value = 18
$ws.Range("R146").Value = 'Hortaliza'

# Row 147
$ws.Range("A147").Value = 12
$ws.Range("B147").Value = 'Mapocho Venta Directa de Santiago'
$ws.Range("C147").Value = 'Metropolitana'
$ws.Range("D147").Value = 44637
$ws.Range("E147").Value = 13
$ws.Range("F147").Value = 100112043
$ws.Range("G147").Value = 'Pepino dulce'
$ws.Range("H147").Value = 'Cultivar IV Región'
$ws.Range("I147").Value = 'Tercera'
$ws.Range("J147").Value = 200
$ws.Range("K147").Value = 8000
$ws.Range("L147").Value = 8000
$ws.Range("M147").Value = 8000
$ws.Range("N147").Value = '$/bandeja 18 kilos'
$ws.Range("O147").Value = 'Provincia de Limarí'
$ws.Range("P147").Value = 444
$ws.Range("Q147").Value = 18
$ws.Range("R147").Value = 'Hortaliza'

Write-Output "done"